$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 116: 2020-08-12, "10:00:00" ---
# Copy the date cell's format (style) from the last existing data row (A115)
# using PasteSpecial(xlPasteFormats) so no new cellXfs entry is created in
# xl/styles.xml - it reuses the already-present style index.
$ws.Range("A115").Copy()
$ws.Range("A116").PasteSpecial(-4122)
$ws.Range("A116").Value = 44055

$ws.Range("B116").Value = "10:00:00"
$ws.Range("C116").Value = 2181
$ws.Range("D116").Value = 1070
$ws.Range("E116").Value = 101
$ws.Range("F116").Value = 2008
$ws.Range("G116").Value = 72

# --- Row 117: 2020-08-14, "10:30:00" ---
$ws.Range("A115").Copy()
$ws.Range("A117").PasteSpecial(-4122)
$ws.Range("A117").Value = 44057

$ws.Range("B117").Value = "10:30:00"
$ws.Range("C117").Value = 2212
$ws.Range("D117").Value = 1088
$ws.Range("E117").Value = 101
$ws.Range("F117").Value = 2028
$ws.Range("G117").Value = 83

$excel.CutCopyMode = $false
